$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate the table header row (A1:E1) from English to German
$ws.Range("A1").Value = "Titel"
$ws.Range("B1").Value = "Sprache"
$ws.Range("C1").Value = "Autor*innen"
$ws.Range("D1").Value = "Lizenz"
$ws.Range("E1").Value = "Veröffentlichungsdatum"

# Remove the (white) solid fill from the data rows, leaving the cell borders intact
$ws.Range("A2:E10").Interior.Pattern = -4142
